$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Rspo3 -> Lgr6 -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rspo3"
$ws.Range("C2").Value = "Lgr6"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.052433
$ws.Range("H2").Value = 0.104866
$ws.Range("I2").Value = 0.02014015227998088
$ws.Range("J2").Value = 0.01351751646668263
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1085923333333333
$ws.Range("N2").Value = 0.325777
$ws.Range("O2").Value = 0.1969429372262032
$ws.Range("P2").Value = 0.2689322520249767
$ws.Range("Q2").Value = 0.005693821813666667
$ws.Range("R2").Value = 0.034162930882
$ws.Range("S2").Value = 0.003966460746202449
$ws.Range("T2").Value = 0.003635296145169666

# Row 3: ECs -> Rspo3 -> Lgr6 -> sCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rspo3"
$ws.Range("C3").Value = "Lgr6"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.052433
$ws.Range("H3").Value = 0.104866
$ws.Range("I3").Value = 0.02014015227998088
$ws.Range("J3").Value = 0.01351751646668263
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.4427975
$ws.Range("N3").Value = 0.885595
$ws.Range("O3").Value = 0.8030570627737967
$ws.Range("P3").Value = 0.7310677479750235
$ws.Range("Q3").Value = 0.0232172013175
$ws.Range("R3").Value = 0.09286880527000001
$ws.Range("S3").Value = 0.01617369153377843
$ws.Range("T3").Value = 0.009882220321512968

# Row 4: FAPs -> Rspo3 -> Lgr6 -> FAPs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo3"
$ws.Range("C4").Value = "Lgr6"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.532784666666667
$ws.Range("H4").Value = 7.598354
$ws.Range("I4").Value = 0.97287335988914
$ws.Range("J4").Value = 0.9794487757202892
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1085923333333333
$ws.Range("N4").Value = 0.325777
$ws.Range("O4").Value = 0.1969429372262032
$ws.Range("P4").Value = 0.2689322520249767
$ws.Range("Q4").Value = 0.2750409967842222
$ws.Range("R4").Value = 2.475368971058
$ws.Range("S4").Value = 0.1916005370456923
$ws.Range("T4").Value = 0.2634053649975637

# Row 5: FAPs -> Rspo3 -> Lgr6 -> sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo3"
$ws.Range("C5").Value = "Lgr6"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.532784666666667
$ws.Range("H5").Value = 7.598354
$ws.Range("I5").Value = 0.97287335988914
$ws.Range("J5").Value = 0.9794487757202892
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4427975
$ws.Range("N5").Value = 0.885595
$ws.Range("O5").Value = 0.8030570627737967
$ws.Range("P5").Value = 0.7310677479750235
$ws.Range("Q5").Value = 1.121510718438333
$ws.Range("R5").Value = 6.729064310630001
$ws.Range("S5").Value = 0.7812728228434477
$ws.Range("T5").Value = 0.7160434107227257

# Row 6: Neutro -> Rspo3 -> Lgr6 -> FAPs
$ws.Range("A6").Value = "Neutro"
$ws.Range("B6").Value = "Rspo3"
$ws.Range("C6").Value = "Lgr6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.01818866666666667
$ws.Range("H6").Value = 0.054566
$ws.Range("I6").Value = 0.006986487830879007
$ws.Range("J6").Value = 0.007033707813028098
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1085923333333333
$ws.Range("N6").Value = 0.325777
$ws.Range("O6").Value = 0.1969429372262032
$ws.Range("P6").Value = 0.2689322520249767
$ws.Range("Q6").Value = 0.001975149753555556
$ws.Range("R6").Value = 0.017776347782
$ws.Range("S6").Value = 0.001375939434308437
$ws.Range("T6").Value = 0.00189159088224332

# Row 7: Neutro -> Rspo3 -> Lgr6 -> sCs
$ws.Range("A7").Value = "Neutro"
$ws.Range("B7").Value = "Rspo3"
$ws.Range("C7").Value = "Lgr6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.01818866666666667
$ws.Range("H7").Value = 0.054566
$ws.Range("I7").Value = 0.006986487830879007
$ws.Range("J7").Value = 0.007033707813028098
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4427975
$ws.Range("N7").Value = 0.885595
$ws.Range("O7").Value = 0.8030570627737967
$ws.Range("P7").Value = 0.7310677479750235
$ws.Range("Q7").Value = 0.008053896128333335
$ws.Range("R7").Value = 0.04832337677
$ws.Range("S7").Value = 0.00561054839657057
$ws.Range("T7").Value = 0.00514211693078478
